$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at 33 (pushes existing rows 33-36 down to 34-37) ---
$ws.Rows("33:33").Insert()

# --- 2. Fill in the new SudachiPy row ---
$ws.Range("A33").Value = "SudachiPy"
$ws.Range("B33").Value = "https://github.com/WorksApplications/sudachi.rs"
$ws.Range("C33").Value = "0.6.2"
$ws.Range("D33").Value = "Works Applications Co., Ltd."
$ws.Range("E33").Value = "Apache-2.0"
$ws.Range("F33").Value = "https://github.com/WorksApplications/sudachi.rs/blob/develop/LICENSE"

# --- 3. Rebuild hyperlinks (row insert does not auto-shift hyperlink anchors) ---
$ws.Hyperlinks.Delete() | Out-Null

$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.crummy.com/software/BeautifulSoup/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/Ousret/charset_normalizer") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "https://github.com/Mimino666/langdetect") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"), "https://github.com/saffsd/langid.py") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B15"), "https://www.numpy.org/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B17"), "https://foss.heptapod.net/openpyxl/openpyxl") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B25"), "https://github.com/python-openxml/python-docx") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B27"), "https://github.com/psf/requests") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B29"), "https://scipy.org/scipylib/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://bazaar.launchpad.net/~leonardr/beautifulsoup/bs4/view/head:/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/Ousret/charset_normalizer/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://github.com/Mimino666/langdetect/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "https://github.com/saffsd/langid.py/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), "https://github.com/numpy/numpy/blob/master/LICENSE.txt") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F17"), "https://foss.heptapod.net/openpyxl/openpyxl/-/blob/branch/3.0/LICENCE.rst") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F25"), "https://github.com/python-openxml/python-docx/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F27"), "https://github.com/requests/requests/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F29"), "https://github.com/scipy/scipy/blob/master/LICENSE.txt") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F24"), "https://docs.python.org/3.8/license.html", "psf-license-agreement-for-python-release") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F19"), "https://github.com/pyinstaller/pyinstaller/blob/develop/COPYING.txt") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F22"), "https://www.riverbankcomputing.com/static/Docs/PyQt5/introduction.html", "license") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B24"), "https://www.python.org/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B19"), "http://www.pyinstaller.org/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B22"), "https://riverbankcomputing.com/software/pyqt/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B11"), "https://matplotlib.org/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B13"), "https://networkx.org/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://matplotlib.org/users/license.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), "https://github.com/networkx/networkx/blob/master/LICENSE.txt") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B37"), "https://github.com/amueller/word_cloud") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F37"), "https://github.com/amueller/word_cloud/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/Esukhia/botok") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/cltk/cltk") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/fxsjy/jieba") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B12"), "https://github.com/taishi-i/nagisa") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B14"), "http://www.nltk.org/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B16"), "https://github.com/yichen0831/opencc-python") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B18"), "https://github.com/lancopku/pkuseg-python") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B21"), "https://pyphen.org/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B20"), "https://github.com/kmike/pymorphy2") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B23"), "https://github.com/PyThaiNLP/pythainlp") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B26"), "https://github.com/natasha/razdel") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B28"), "https://github.com/alvations/sacremoses") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B30"), "https://spacy.io/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B31"), "https://github.com/ponrawee/ssg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B34"), "https://github.com/sloria/TextBlob") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B35"), "https://github.com/mideind/Tokenizer") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B36"), "https://github.com/undertheseanlp/underthesea") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/Esukhia/botok/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://github.com/cltk/cltk/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://github.com/fxsjy/jieba/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://github.com/taishi-i/nagisa/blob/master/LICENSE.txt") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), "https://github.com/nltk/nltk/blob/develop/LICENSE.txt") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), "https://github.com/yichen0831/opencc-python/blob/master/LICENSE.txt") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F18"), "https://github.com/lancopku/pkuseg-python/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F21"), "https://github.com/Kozea/Pyphen/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F20"), "https://github.com/kmike/pymorphy2/", "pymorphy2") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F23"), "https://github.com/PyThaiNLP/pythainlp/blob/dev/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F26"), "https://github.com/natasha/razdel/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F28"), "https://github.com/alvations/sacremoses/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F30"), "https://github.com/explosion/spaCy/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F31"), "https://github.com/ponrawee/ssg/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F34"), "https://github.com/sloria/TextBlob/blob/dev/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F35"), "https://github.com/mideind/Tokenizer/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F36"), "https://github.com/undertheseanlp/underthesea/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/Xangis/extra-stopwords") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B10"), "https://github.com/michmech/lemmatization-lists") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B32"), "https://github.com/stopwords-iso/stopwords-iso") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://github.com/Xangis/extra-stopwords/blob/master/LICENSE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://github.com/michmech/lemmatization-lists/blob/master/LICENCE") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F32"), "https://github.com/stopwords-iso/stopwords-iso/blob/master/LICENSE") | Out-Null

# --- 4. Update the sort state to cover the new row range ---
$sortObj = $ws.Sort
$sortObj.SetRange($ws.Range("A2:F37"))
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A37")) | Out-Null
$sortObj.Header = 2
$sortObj.Apply()

# --- 5. Update the active selection to C33 ---
$ws.Range("C33").Select()
